# VET - Updated Monte Carlo simulations
# Delete the "E_BIO" / "BIO" row (row 6) from the power-plant related sheets.
$wb = $excel.ActiveWorkbook

$sheetNames = @("PowerPlants", "PowerPlantsPerformance", "PowerPlantsCosts", "PowerPlantsConstraints")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Rows.Item(6).Delete()
}
